$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values for row 2
$ws.Range("A2").Value = 0
$ws.Range("C2").Value = 1
$ws.Range("F2").Value = 43831
$ws.Range("G2").Value = 43834

# Center the date cells vertically (in addition to existing horizontal centering)
$ws.Range("F2:G2").VerticalAlignment = -4108

# Update the active selection shown in the sheet view
$ws.Range("C7").Select() | Out-Null
